$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete iAuthor test-case rows (old rows 7-13); only the
# header plus five data rows remain (TC_205, TC_180, TC_205, TC_206).
$ws.Range("A7:C13").EntireRow.Delete()

# Row 2: Pre-Request Login/Create Exam now maps to iAU_TC_ID_205
$ws.Range("A2").Value = " iAU_TC_ID_205"

# Row 3: Pre-Request Add New Users now maps to iAU_TC_ID_205
$ws.Range("A3").Value = " iAU_TC_ID_205"

# Row 4: newly added case - iAU_TC_ID_180 / Edit user
$ws.Range("A4").Value = "iAU_TC_ID_180"
$ws.Range("B4").Value = "@RegressionA Validation of Manage Delivery --> Edit user "

# Row 5: newly added case - iAU_TC_ID_205 / Venue Summary
$ws.Range("A5").Value = "iAU_TC_ID_205"
$ws.Range("B5").Value = "@RegressionA Validation of Delivery --> Venue Summary "

# Row 6: newly added case - iAU_TC_ID_206 / Live Monitor Dashboard
$ws.Range("A6").Value = "iAU_TC_ID_206"
$ws.Range("B6").Value = "@RegressionA Validation of Delivery --> Live Monitor Dashboard "
